# "added some more pre-processing"
#
# 1. Rename the existing "results" sheet to "results_allgenres".
# 2. Add a new blank worksheet named "Feuil3" right after it.
# 3. On "results_allgenres", scroll back to the top-left (drop the
#    stashed topLeftCell="H2") and move the active selection to J5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the first/active sheet.
$ws.Name = "results_allgenres"

# 2. Insert a brand-new empty sheet immediately after it, named "Feuil3".
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$newSheet.Name = "Feuil3"

# 3. Re-activate the results sheet and move/scroll the selection to J5
#    (this also resets any stashed topLeftCell scroll position).
$ws.Activate()
$ws.Range("J5").Select()
